$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for account "005009922" (ANA, balance 6538.79) was removed from
# the export, so delete the whole row here too; every row below it shifts
# up by one to fill the gap.
$target = $ws.Cells.Find("005009922")
if ($target -ne $null) {
    $target.EntireRow.Delete()
}
